# Auto-generated edit script: updates Leve profit-calculation values
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns) across
# all eight crafting-job sheets, per refreshed market-board data.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 249.66667
$ws.Range("I12").Value = 249.66667
$ws.Range("K12").Value = 249.66667
$ws.Range("M12").Value = -79.66667000000001
$ws.Range("H43").Value = 899
$ws.Range("J43").Value = 999.5
$ws.Range("L43").Value = 999.5
$ws.Range("N43").Value = -1137.5
$ws.Range("H53").Value = 262.16666
$ws.Range("I53").Value = 218.8
$ws.Range("J53").Value = 293.14285
$ws.Range("K53").Value = 218.8
$ws.Range("L53").Value = 293.14285
$ws.Range("M53").Value = 418.2
$ws.Range("N53").Value = -1567.14285
$ws.Range("H55").Value = 764.3333
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 764.3333
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 764.3333
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -1192.3333
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880
$ws.Range("H69").Value = 107998.6
$ws.Range("I69").Value = 5552.5
$ws.Range("K69").Value = 16657.5
$ws.Range("M69").Value = -15783.5
$ws.Range("H72").Value = 107998.6
$ws.Range("I72").Value = 5552.5
$ws.Range("K72").Value = 49972.5
$ws.Range("M72").Value = -45604.5
$ws.Range("H88").Value = 9000
$ws.Range("J88").Value = 8000
$ws.Range("L88").Value = 8000
$ws.Range("N88").Value = -8812
$ws.Range("H91").Value = 9000
$ws.Range("J91").Value = 8000
$ws.Range("L91").Value = 8000
$ws.Range("N91").Value = -10808
$ws.Range("H92").Value = 798.65
$ws.Range("J92").Value = 603.4286
$ws.Range("L92").Value = 603.4286
$ws.Range("N92").Value = -3099.4286
$ws.Range("H94").Value = 1702.8572
$ws.Range("I94").Value = 1530
$ws.Range("K94").Value = 1530
$ws.Range("M94").Value = -1079
$ws.Range("H98").Value = 4205.6665
$ws.Range("J98").Value = 10500
$ws.Range("L98").Value = 10500
$ws.Range("N98").Value = -13496
$ws.Range("H104").Value = 491
$ws.Range("I104").Value = 491
$ws.Range("K104").Value = 1473
$ws.Range("M104").Value = 274
$ws.Range("H106").Value = 4999.5
$ws.Range("I106").Value = 4999
$ws.Range("K106").Value = 4999
$ws.Range("M106").Value = -4368
$ws.Range("H107").Value = 3033.6
$ws.Range("I107").Value = 6365.5
$ws.Range("J107").Value = 812.3333
$ws.Range("K107").Value = 6365.5
$ws.Range("L107").Value = 812.3333
$ws.Range("M107").Value = -4445.5
$ws.Range("N107").Value = -4652.3333
$ws.Range("H113").Value = 7914.2856
$ws.Range("J113").Value = 7250
$ws.Range("L113").Value = 7250
$ws.Range("N113").Value = -13758
$ws.Range("H122").Value = 4205.6665
$ws.Range("J122").Value = 10500
$ws.Range("L122").Value = 31500
$ws.Range("N122").Value = -36400

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7057.0835
$ws.Range("I74").Value = 6789.636
$ws.Range("K74").Value = 6789.636
$ws.Range("M74").Value = -5915.636
$ws.Range("H77").Value = 7057.0835
$ws.Range("I77").Value = 6789.636
$ws.Range("K77").Value = 33948.18
$ws.Range("M77").Value = -29580.18
$ws.Range("H88").Value = 3500
$ws.Range("I88").Value = 2333.3333
$ws.Range("J88").Value = 4666.6665
$ws.Range("K88").Value = 2333.3333
$ws.Range("L88").Value = 4666.6665
$ws.Range("M88").Value = -1927.3333
$ws.Range("N88").Value = -5478.6665
$ws.Range("H91").Value = 3500
$ws.Range("I91").Value = 2333.3333
$ws.Range("J91").Value = 4666.6665
$ws.Range("K91").Value = 2333.3333
$ws.Range("L91").Value = 4666.6665
$ws.Range("M91").Value = -929.3332999999998
$ws.Range("N91").Value = -7474.6665
$ws.Range("H102").Value = 7825.857
$ws.Range("I102").Value = 8296.833000000001
$ws.Range("K102").Value = 8296.833000000001
$ws.Range("M102").Value = -6674.833000000001
$ws.Range("H122").Value = 3357.8
$ws.Range("I122").Value = 3215.3635
$ws.Range("J122").Value = 3749.5
$ws.Range("K122").Value = 9646.0905
$ws.Range("L122").Value = 11248.5
$ws.Range("M122").Value = -7196.0905
$ws.Range("N122").Value = -16148.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 731.3333
$ws.Range("I5").Value = 731.3333
$ws.Range("K5").Value = 731.3333
$ws.Range("M5").Value = -618.3333
$ws.Range("H7").Value = 17550
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 35000
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 35000
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -35226
$ws.Range("H20").Value = 4137.4614
$ws.Range("I20").Value = 1283
$ws.Range("J20").Value = 6584.143
$ws.Range("K20").Value = 1283
$ws.Range("L20").Value = 6584.143
$ws.Range("M20").Value = -1036
$ws.Range("N20").Value = -7078.143
$ws.Range("H86").Value = 6966.5
$ws.Range("I86").Value = 2699.75
$ws.Range("J86").Value = 15500
$ws.Range("K86").Value = 2699.75
$ws.Range("L86").Value = 15500
$ws.Range("M86").Value = -1576.75
$ws.Range("N86").Value = -17746
$ws.Range("H89").Value = 6966.5
$ws.Range("I89").Value = 2699.75
$ws.Range("J89").Value = 15500
$ws.Range("K89").Value = 13498.75
$ws.Range("L89").Value = 77500
$ws.Range("M89").Value = -7882.75
$ws.Range("N89").Value = -88732
$ws.Range("H94").Value = 3420.9
$ws.Range("I94").Value = 2368.5
$ws.Range("J94").Value = 4999.5
$ws.Range("K94").Value = 2368.5
$ws.Range("L94").Value = 4999.5
$ws.Range("M94").Value = -1917.5
$ws.Range("N94").Value = -5901.5
$ws.Range("H105").Value = 3356
$ws.Range("I105").Value = 3264
$ws.Range("K105").Value = 3264
$ws.Range("M105").Value = -1517
$ws.Range("H107").Value = 1245
$ws.Range("I107").Value = 1245
$ws.Range("K107").Value = 1245
$ws.Range("M107").Value = 675

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 64.666664
$ws.Range("J7").Value = 14
$ws.Range("L7").Value = 14
$ws.Range("N7").Value = -240
$ws.Range("H16").Value = 842.1429000000001
$ws.Range("I16").Value = 908
$ws.Range("K16").Value = 908
$ws.Range("M16").Value = -621
$ws.Range("H99").Value = 7000
$ws.Range("I99").Value = 7000
$ws.Range("K99").Value = 7000
$ws.Range("M99").Value = -5502
$ws.Range("H105").Value = 4098.7
$ws.Range("I105").Value = 4338.4
$ws.Range("K105").Value = 4338.4
$ws.Range("M105").Value = -2591.4
$ws.Range("H113").Value = 842.1429000000001
$ws.Range("I113").Value = 908
$ws.Range("K113").Value = 908
$ws.Range("M113").Value = 1262
$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -18530

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3136
$ws.Range("J122").Value = 3476.0557
$ws.Range("L122").Value = 31284.5013
$ws.Range("N122").Value = -36184.5013
$ws.Range("H131").Value = 3109.3333
$ws.Range("J131").Value = 3109.3333
$ws.Range("L131").Value = 9327.999899999999
$ws.Range("N131").Value = -19407.9999
$ws.Range("H140").Value = 2546
$ws.Range("I140").Value = 2546
$ws.Range("K140").Value = 7638
$ws.Range("M140").Value = -2458

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 10000
$ws.Range("I97").Value = 10000
$ws.Range("K97").Value = 10000
$ws.Range("M97").Value = -9504

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6081.4707
$ws.Range("I22").Value = 4648.3335
$ws.Range("K22").Value = 4648.3335
$ws.Range("M22").Value = -4353.3335
$ws.Range("H27").Value = 6081.4707
$ws.Range("I27").Value = 4648.3335
$ws.Range("K27").Value = 4648.3335
$ws.Range("M27").Value = -4541.3335
$ws.Range("H40").Value = 3299.125
$ws.Range("I40").Value = 3299.125
$ws.Range("K40").Value = 3299.125
$ws.Range("M40").Value = -3163.125
$ws.Range("H55").Value = 2235.1538
$ws.Range("I55").Value = 5321
$ws.Range("K55").Value = 5321
$ws.Range("M55").Value = -5148
$ws.Range("H132").Value = 4539
$ws.Range("I132").Value = 3390.8
$ws.Range("K132").Value = 10172.4
$ws.Range("M132").Value = -7642.400000000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1024.5
$ws.Range("I113").Value = 1112.8889
$ws.Range("K113").Value = 3338.6667
$ws.Range("M113").Value = -1168.6667
$ws.Range("H126").Value = 1496.8695
$ws.Range("I126").Value = 1366.4546
$ws.Range("J126").Value = 1616.4166
$ws.Range("K126").Value = 4099.3638
$ws.Range("L126").Value = 4849.2498
$ws.Range("M126").Value = -1629.3638
$ws.Range("N126").Value = -9789.2498

Write-Host "Applied 236 cell updates"
